# Updates cryptocurrency price (D) and 1h volume change (E) columns
# to reflect latest scraped values, as produced by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.407.09"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "'1.619.00"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("D5").Value = "'212.86"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D10").Value = "'19.18"
$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("E11").Value = "  -0.66%  "

$ws.Range("D12").Value = "'1.845.99"
$ws.Range("E12").Value = "  +1.10%  "

$ws.Range("D13").Value = "'1.620.29"
$ws.Range("E13").Value = "  +1.43%  "

$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").Value = "'0.509"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").Value = "'63.83"
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("D17").Value = "'236.17"
$ws.Range("E17").Value = "  +8.91%  "

$ws.Range("D18").Value = "'26.408.72"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("E19").Value = "  +4.20%  "

$ws.Range("D20").Value = "'0.0₃0726"
$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("E22").Value = "  -0.66%  "

$ws.Range("D23").Value = "'9.12"
$ws.Range("E23").Value = "  +1.11%  "

$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  +2.84%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "'7.06"
$ws.Range("E27").Value = "  +0.75%  "

$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").Value = "'15.56"
$ws.Range("E29").Value = "  +2.27%  "

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  -0.50%  "

$ws.Range("D32").Value = "'1.515.46"
$ws.Range("E32").Value = "  +5.59%  "

$ws.Range("E33").Value = "  +1.28%  "

$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("E35").Value = "  +2.90%  "

$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("E37").Value = "  +1.99%  "

$ws.Range("E38").Value = "  +0.31%  "

$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("D40").Value = "'5.89"
$ws.Range("E40").Value = "  +1.69%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("D43").Value = "'1.758.38"
$ws.Range("E43").Value = "  +1.35%  "

$ws.Range("D44").Value = "'0.762"
$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("D45").Value = "'62.00"
$ws.Range("E45").Value = "  +1.35%  "

$ws.Range("D46").Value = "'0.908"
$ws.Range("E46").Value = "  +0.41%  "

$ws.Range("D47").Value = "'90.26"
$ws.Range("E47").Value = "  +3.54%  "

$ws.Range("E48").Value = "  +1.66%  "

$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").Value = "'0.0965"
$ws.Range("E50").Value = "  +0.88%  "

$ws.Range("D51").Value = "'7.51"
$ws.Range("E51").Value = "  +0.82%  "
